# Auto-generated Excel COM-interop script
# Applies updated "想去人数" (F column) counts across all 4 worksheets
$wb = $excel.ActiveWorkbook

# --- Worksheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 261
$ws.Range("F3").Value = 611
$ws.Range("F5").Value = 173
$ws.Range("F6").Value = 2807
$ws.Range("F8").Value = 52
$ws.Range("F9").Value = 28
$ws.Range("F10").Value = 371
$ws.Range("F11").Value = 20
$ws.Range("F12").Value = 309
$ws.Range("F14").Value = 5898
$ws.Range("F15").Value = 628
$ws.Range("F16").Value = 1042
$ws.Range("F17").Value = 4
$ws.Range("F18").Value = 96
$ws.Range("F19").Value = 170
$ws.Range("F21").Value = 519
$ws.Range("F22").Value = 17
$ws.Range("F23").Value = 18
$ws.Range("F24").Value = 7
$ws.Range("F25").Value = 1292
$ws.Range("F27").Value = 2
$ws.Range("F28").Value = 27
$ws.Range("F29").Value = 2045
$ws.Range("F30").Value = 165
$ws.Range("F33").Value = 3259

# --- Worksheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 342
$ws.Range("F8").Value = 78
$ws.Range("F12").Value = 641
$ws.Range("F17").Value = 61
$ws.Range("F18").Value = 625
$ws.Range("F20").Value = 60
$ws.Range("F24").Value = 4045
$ws.Range("F26").Value = 12
$ws.Range("F28").Value = 117
$ws.Range("F30").Value = 67
$ws.Range("F33").Value = 16

# --- Worksheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F8").Value = 1477
$ws.Range("F12").Value = 623

# --- Worksheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F7").Value = 1477
$ws.Range("F11").Value = 261
$ws.Range("F12").Value = 611
$ws.Range("F13").Value = 2807
$ws.Range("F14").Value = 52
$ws.Range("F15").Value = 28
$ws.Range("F16").Value = 623
$ws.Range("F17").Value = 371
$ws.Range("F18").Value = 78
$ws.Range("F19").Value = 20
$ws.Range("F20").Value = 309
$ws.Range("F22").Value = 5899
$ws.Range("F23").Value = 628
$ws.Range("F24").Value = 1042
$ws.Range("F25").Value = 96
$ws.Range("F26").Value = 170
$ws.Range("F28").Value = 519
$ws.Range("F30").Value = 61
$ws.Range("F32").Value = 60
$ws.Range("F33").Value = 17
$ws.Range("F36").Value = 1292
$ws.Range("F37").Value = 12
$ws.Range("F39").Value = 117
$ws.Range("F41").Value = 27
$ws.Range("F42").Value = 67
$ws.Range("F44").Value = 2045
$ws.Range("F47").Value = 165
$ws.Range("F50").Value = 3259
